$d = $word.ActiveDocument

# Find the paragraph that begins "Since 2008, Prof. Lawson..." -- this is
# the opening paragraph of the (as yet untitled) Law & Literature section.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Since 2008, Prof. Lawson*") {
        $target = $p
        break
    }
}
$targetIndex = $target.Index

# Insert a new empty paragraph immediately before it. After this call, the
# original paragraph object/position becomes the new (empty) paragraph, and
# the "Since 2008..." text is pushed down to the next paragraph index.
$target.Range.InsertParagraphBefore()

$heading = $target
$heading.Style = "Heading 2"
$heading.Range.Text = "Law & Literature"

# Add a (zero-length) bookmark at the very start of the new heading,
# matching the bookmark convention used for the other Heading2/Heading3
# section headings in this document.
$bmRange = $d.Range($heading.Range.Start, $heading.Range.Start)
$d.Bookmarks.Add("law-literature", $bmRange)

# The paragraph that used to open the section (styled as BodyText) now
# acts as the section's first paragraph, so restyle it to match the other
# sections' opening paragraphs.
$body = $d.Paragraphs.Item($targetIndex + 1)
$body.Style = "First Paragraph"
